# Ajout ref connecteur alim
# "Connecteur d'alimentation" row (row 12): replace the placeholder
# "Diametre a definir" in REF 2 (G12) with the actual connector reference,
# formatted as text, and fill in the matching Prix 2 CHF (H12).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$refCell = $ws.Range("G12")
$refCell.NumberFormat = "@"
$refCell.Value = "694106301002"

$ws.Range("H12").Value = 0.921
